$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text assignments (URLs, coin names, and price strings that Excel cannot
#     mis-parse as numbers, e.g. thousands-grouped "63.745.25") ---
$ws.Range("D2").Value = "63.745.25"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").Value = "3.481.48"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("D7").Value = "3.480.18"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").Value = "4.066.56"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Value = "3.484.23"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "63.734.06"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "3.617.90"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").Value = "3.480.92"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").Value = "2.432.08"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -0.80%  "

# --- Price strings that look like plain numbers (e.g. "0.998", "7.17").
#     These must stay text cells (t="inlineStr" in the source file), so each one is
#     switched to a Text number format right before the write (preventing COM from
#     coercing it into a floating point value) and switched back to the default
#     "Normal" style immediately after, so no visible/style change is left behind. ---
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "580.63"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "130.32"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.488"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.123"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "7.17"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.380"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "27.30"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0000177"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "10.06"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "14.32"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.63"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "383.05"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.574"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "72.81"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.0000111"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.57"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.39"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "8.16"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "23.60"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.143"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.24"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.55"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "167.51"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.85"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0796"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "26.84"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.809"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.20"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "41.19"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.36"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.63"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "6.82"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.886"
$cell.Style = "Normal"
